$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030043103348618
$ws.Range("D2").Value = 1.034060109576233
$ws.Range("E2").Value = 1.033679356161622
$ws.Range("F2").Value = 1.03945324294509
$ws.Range("I2").Value = 1.033898326160537
$ws.Range("J2").Value = 1.035186531408275
$ws.Range("K2").Value = 1.036860669823648
$ws.Range("L2").Value = 1.036481011334414
$ws.Range("M2").Value = 1.042238387803358
$ws.Range("N2").Value = 1.015472065551848
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031560454322792
$ws.Range("D3").Value = 1.035227976904537
$ws.Range("E3").Value = 1.035145621879192
$ws.Range("F3").Value = 1.041094975434854
$ws.Range("I3").Value = 1.03432305468253
$ws.Range("J3").Value = 1.036342271214647
$ws.Range("K3").Value = 1.037836690274446
$ws.Range("L3").Value = 1.037754554702258
$ws.Range("M3").Value = 1.043688151422883
$ws.Range("N3").Value = 1.015878601457215
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032540104785392
$ws.Range("D4").Value = 1.03598159183038
$ws.Range("E4").Value = 1.036092612784714
$ws.Range("F4").Value = 1.042155349343132
$ws.Range("I4").Value = 1.034595265467389
$ws.Range("J4").Value = 1.037087598825487
$ws.Range("K4").Value = 1.038465603081561
$ws.Range("L4").Value = 1.038576343296949
$ws.Range("M4").Value = 1.044623846811822
$ws.Range("N4").Value = 1.016140200613545
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032951438602859
$ws.Range("D5").Value = 1.036297921841922
$ws.Range("E5").Value = 1.036490310233573
$ws.Range("F5").Value = 1.042600677555253
$ws.Range("I5").Value = 1.03470907959447
$ws.Range("J5").Value = 1.037400340421855
$ws.Range("K5").Value = 1.038729372801183
$ws.Range("L5").Value = 1.03892128567925
$ws.Range("M5").Value = 1.045016648604854
$ws.Range("N5").Value = 1.016249830834757
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033020473669664
$ws.Range("D6").Value = 1.036351006518666
$ws.Range("E6").Value = 1.036557061157427
$ws.Range("F6").Value = 1.042675423867946
$ws.Range("I6").Value = 1.034728153030359
$ws.Range("J6").Value = 1.037452816490562
$ws.Range("K6").Value = 1.03877362441469
$ws.Range("L6").Value = 1.038979171735025
$ws.Range("M6").Value = 1.045082568928962
$ws.Range("N6").Value = 1.016268218015649
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032545603044523
$ws.Range("D7").Value = 1.035985820563491
$ws.Range("E7").Value = 1.036097928467162
$ws.Range("F7").Value = 1.042161301607215
$ws.Range("I7").Value = 1.034596788702476
$ws.Range("J7").Value = 1.037091780019224
$ws.Range("K7").Value = 1.038469130034502
$ws.Range("L7").Value = 1.038580954536164
$ws.Range("M7").Value = 1.044629097654971
$ws.Range("N7").Value = 1.016141666854403
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030556356009267
$ws.Range("D8").Value = 1.034455229127037
$ws.Range("E8").Value = 1.034175261659002
$ws.Range("F8").Value = 1.040008481610951
$ws.Range("I8").Value = 1.034042408961948
$ws.Range("J8").Value = 1.035577644057593
$ws.Range("K8").Value = 1.037191070429418
$ws.Range("L8").Value = 1.036911888069478
$ws.Range("M8").Value = 1.042728844706125
$ws.Range("N8").Value = 1.01560975951882
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027033926408934
$ws.Range("D9").Value = 1.031741952004229
$ws.Range("E9").Value = 1.030773237266643
$ws.Range("F9").Value = 1.036199622837501
$ws.Range("I9").Value = 1.033045332614541
$ws.Range("J9").Value = 1.032889953196096
$ws.Range("K9").Value = 1.034918482936103
$ws.Range("L9").Value = 1.033952960222026
$ws.Range("M9").Value = 1.039361549557324
$ws.Range("N9").Value = 1.014661192793068
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024673498351149
$ws.Range("D10").Value = 1.029921799288834
$ws.Range("E10").Value = 1.028495220527489
$ws.Range("F10").Value = 1.033649410710062
$ws.Range("I10").Value = 1.032366835855045
$ws.Range("J10").Value = 1.031084508906606
$ws.Range("K10").Value = 1.033389256056277
$ws.Range("L10").Value = 1.031967848844481
$ws.Range("M10").Value = 1.037103428738991
$ws.Range("N10").Value = 1.014021061959145
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023648386189367
$ws.Range("D11").Value = 1.029130877410357
$ws.Range("E11").Value = 1.027506316492274
$ws.Range("F11").Value = 1.032542386772945
$ws.Range("I11").Value = 1.032069724625474
$ws.Range("J11").Value = 1.03029938644323
$ws.Range("K11").Value = 1.032723634150022
$ws.Range("L11").Value = 1.031105195811052
$ws.Range("M11").Value = 1.036122355782898
$ws.Range("N11").Value = 1.013742001100493
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023267146251724
$ws.Range("D12").Value = 1.028836667280061
$ws.Range("E12").Value = 1.027138605824287
$ws.Range("F12").Value = 1.032130760220789
$ws.Range("I12").Value = 1.031958861764741
$ws.Range("J12").Value = 1.030007243872232
$ws.Range("K12").Value = 1.032475865437336
$ws.Range("L12").Value = 1.030784294340017
$ws.Range("M12").Value = 1.035757435139102
$ws.Range("N12").Value = 1.013638059685872
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023348944871667
$ws.Range("D13").Value = 1.028899795787308
$ws.Range("E13").Value = 1.02721749873046
$ws.Range("F13").Value = 1.032219075119354
$ws.Range("I13").Value = 1.031982665030853
$ws.Range("J13").Value = 1.030069932816822
$ws.Range("K13").Value = 1.032529036645576
$ws.Range("L13").Value = 1.030853150393233
$ws.Range("M13").Value = 1.035835734967194
$ws.Range("N13").Value = 1.013660368458553
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023616882388014
$ws.Range("D14").Value = 1.029106566658613
$ws.Range("E14").Value = 1.027475929385799
$ws.Range("F14").Value = 1.032508370412839
$ws.Range("I14").Value = 1.032060570941874
$ws.Range("J14").Value = 1.030275248369326
$ws.Range("K14").Value = 1.032703164308965
$ws.Range("L14").Value = 1.031078679717284
$ws.Range("M14").Value = 1.036092201729756
$ws.Range("N14").Value = 1.013733415120128
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023781905200534
$ws.Range("D15").Value = 1.029233908247605
$ws.Range("E15").Value = 1.027635105342299
$ws.Range("F15").Value = 1.032686557606419
$ws.Range("I15").Value = 1.032108504633062
$ws.Range("J15").Value = 1.030401681781289
$ws.Range("K15").Value = 1.03281038002906
$ws.Range("L15").Value = 1.031217572703808
$ws.Range("M15").Value = 1.036250151896454
$ws.Range("N15").Value = 1.013778383609219
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024741466763418
$ws.Range("D16").Value = 1.029974230812133
$ws.Range("E16").Value = 1.028560797066558
$ws.Range("F16").Value = 1.033722820864877
$ws.Range("I16").Value = 1.032386483902372
$ws.Range("J16").Value = 1.031136543454952
$ws.Range("K16").Value = 1.03343335770103
$ws.Range("L16").Value = 1.032025034402085
$ws.Range("M16").Value = 1.037168468952153
$ws.Range("N16").Value = 1.014039542415781
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025342554235766
$ws.Range("D17").Value = 1.030437864476788
$ws.Range("E17").Value = 1.029140779957857
$ws.Range("F17").Value = 1.034372091866654
$ws.Range("I17").Value = 1.032559961856622
$ws.Range("J17").Value = 1.031596598453606
$ws.Range("K17").Value = 1.033823204239409
$ws.Range("L17").Value = 1.032530700376939
$ws.Range("M17").Value = 1.037743615309382
$ws.Range("N17").Value = 1.014202854762011
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02569286675069
$ws.Range("D18").Value = 1.030708026124561
$ws.Range("E18").Value = 1.029478832947077
$ws.Range("F18").Value = 1.034750534743204
$ws.Range("I18").Value = 1.03266082878505
$ws.Range("J18").Value = 1.031864617855303
$ws.Range("K18").Value = 1.034050262150861
$ws.Range("L18").Value = 1.03282534964204
$ws.Range("M18").Value = 1.038078771655011
$ws.Range("N18").Value = 1.014297930926823
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025812265166568
$ws.Range("D19").Value = 1.030800099052842
$ws.Range("E19").Value = 1.029594059693928
$ws.Range("F19").Value = 1.034879529275592
$ws.Range("I19").Value = 1.032695167675392
$ws.Range("J19").Value = 1.0319559510578
$ws.Range("K19").Value = 1.034127626786198
$ws.Range("L19").Value = 1.032925767326775
$ws.Range("M19").Value = 1.038192997932591
$ws.Range("N19").Value = 1.014330318814028
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02527809349254
$ws.Range("D20").Value = 1.030388148755477
$ws.Range("E20").Value = 1.029078578311866
$ws.Range("F20").Value = 1.034302458832412
$ws.Range("I20").Value = 1.032541382424659
$ws.Range("J20").Value = 1.031547272372416
$ws.Range("K20").Value = 1.033781411893189
$ws.Range("L20").Value = 1.032476478010482
$ws.Range("M20").Value = 1.037681940402027
$ws.Range("N20").Value = 1.014185351655622
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023537994475494
$ws.Range("D21").Value = 1.029045689632048
$ws.Range("E21").Value = 1.027399838789037
$ws.Range("F21").Value = 1.032423192096106
$ws.Range("I21").Value = 1.03203764349788
$ws.Range("J21").Value = 1.030214802296969
$ws.Range("K21").Value = 1.032651902665609
$ws.Range("L21").Value = 1.031012280110298
$ws.Range("M21").Value = 1.036016692737681
$ws.Range("N21").Value = 1.013711912615023
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022441213409959
$ws.Range("D22").Value = 1.02819916035781
$ws.Range("E22").Value = 1.026342100150132
$ws.Range("F22").Value = 1.031239136334464
$ws.Range("I22").Value = 1.031718013353716
$ws.Range("J22").Value = 1.029374051959243
$ws.Range("K22").Value = 1.031938681829184
$ws.Range("L22").Value = 1.030088934935417
$ws.Range("M22").Value = 1.034966749194632
$ws.Range("N22").Value = 1.01341258754947
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02302289895955
$ws.Range("D23").Value = 1.028648158648974
$ws.Range("E23").Value = 1.026903044216543
$ws.Range("F23").Value = 1.031867066771559
$ws.Range("I23").Value = 1.03188773250074
$ws.Range("J23").Value = 1.02982003462237
$ws.Range("K23").Value = 1.032317065752057
$ws.Range("L23").Value = 1.030578681314485
$ws.Range("M23").Value = 1.035523626463525
$ws.Range("N23").Value = 1.013571423442333
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025307221428903
$ws.Range("D24").Value = 1.030410613978259
$ws.Range("E24").Value = 1.029106685303691
$ws.Range("F24").Value = 1.034333923824953
$ws.Range("I24").Value = 1.032549778658515
$ws.Range("J24").Value = 1.031569561703078
$ws.Range("K24").Value = 1.033800297084483
$ws.Range("L24").Value = 1.032500979682135
$ws.Range("M24").Value = 1.037709809616799
$ws.Range("N24").Value = 1.014193261115814
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027946652436432
$ws.Range("D25").Value = 1.032445359818908
$ws.Range("E25").Value = 1.031654463336633
$ws.Range("F25").Value = 1.037186188244094
$ws.Range("I25").Value = 1.033305514207381
$ws.Range("J25").Value = 1.033587156266984
$ws.Range("K25").Value = 1.035508469239016
$ws.Range("L25").Value = 1.034720078111733
$ws.Range("M25").Value = 1.040234367495044
$ws.Range("N25").Value = 1.014907773973715
